$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 2
    7  = 1
    8  = 3
    9  = 3
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 2
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 0
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 3
    28 = 1
    29 = 0
    30 = 2
    31 = 1
    32 = 1
    33 = 1
    34 = 0
    36 = 2
    37 = 3
    38 = 0
    39 = 1
    40 = 0
    41 = 1
    42 = 1
    43 = 2
    44 = 1
    45 = 1
    46 = 0
    47 = 2
    48 = 0
    49 = 4
    50 = 0
    51 = 1
    53 = 1
    54 = 2
    55 = 1
    56 = 2
    57 = 2
    58 = 1
    59 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
